$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before Z, shifting old Z->AA and old AA->AB
$ws.Columns("Z:Z").Insert()

# New header for the inserted column
$ws.Range("Z1").Value = "MgCa Coretop modelled temperature"

# New value for the inserted data cell
$ws.Range("Z2").Value = 27.5334

# Updated values in R2:T2
$ws.Range("R2").Value = 27.59
$ws.Range("S2").Value = -3.07169423421227
$ws.Range("T2").Value = -3.54476090087887

# Updated values for the shifted columns (previously Z2/AA2)
$ws.Range("AA2").Value = -3.020066670000002
$ws.Range("AB2").Value = -3.493133329999999
